$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-removed rows (original rows 240-251) first.
$ws.Rows("240:251").Delete()

# Update rows 229-239 with the new set of URLs (replacing the old tail of the list).
$ws.Range("A229").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/cankareva/63b2ee9fb65028de5608665f"
$ws.Range("A230").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/liman-4-novi-sad/63aec89f2562e3893a016eb1"
$ws.Range("A231").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/63aec0a83c109c18d20bbe5a"
$ws.Range("A232").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/nova-detelinara/63aeb4ae568241ab88032ed5"
$ws.Range("A233").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/bulevar-oslobodjenja/63ac44313036173803099566"
$ws.Range("A234").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/centar-novi-sad/63ac7a082af8a275fc0fbaf8"
$ws.Range("A235").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/gunduliceva/63a45a77d20f9e794102a92a"
$ws.Range("A236").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/novo-naselje/63aedb599b08bc373508ac6a"
$ws.Range("A237").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/telep/63ad8a1e40c4f3781a04c17e"
$ws.Range("A238").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/mise-dimitrijevica/6339fa175936bf895b08c5c1"
$ws.Range("A239").Value = "https://www.4zida.rs/izdavanje/stanovi/novi-sad/oglas/socijalno/63ac37f0413a0420b701da8c"

# Restore the active selection to match the saved view state.
$ws.Range("L243").Select()
